# Weekly data refresh: a new week's price observation is inserted as the
# (new) row 40, pushing the existing rows 40-56 down to 41-57.
#
# Result: dimension grows from A1:T56 to A1:T57, and the row that used to
# be the last one (old row 56, "Tercera") becomes the new last row (57).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 40 - this shifts rows 40..56 down to
# 41..57 and grows the sheet's used range accordingly.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly record.
$ws.Range("A40").Value = 6
$ws.Range("B40").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 44627
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100104
$ws.Range("H40").Value = "Frutos de pepita"
$ws.Range("I40").Value = 100104003
$ws.Range("J40").Value = "Membrillo"
$ws.Range("K40").Value = "Champion"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 12
$ws.Range("N40").Value = 350000
$ws.Range("O40").Value = 350000
$ws.Range("P40").Value = 350000
$ws.Range("Q40").Value = "`$/bins (450 kilos)"
$ws.Range("R40").Value = "Región de O'Higgins"
$ws.Range("S40").Value = 778
$ws.Range("T40").Value = 450
